# Daily attendance processing - 2025-10-18 21:16:52
# For each row in the "Recorded By" column (G), when the value begins with
# "System, " (i.e. "System" is the first comma-separated entry), move that
# first entry so it becomes the second entry instead - swapping it with the
# entry that currently follows it. Any further entries keep their relative
# order and position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)  # Column G = 7
    $value = $cell.Value2

    if ($null -eq $value) { continue }
    if ($value -isnot [string]) { continue }

    if ($value.StartsWith("System, ")) {
        $parts = $value.Split(",")
        for ($i = 0; $i -lt $parts.Length; $i++) {
            $parts[$i] = $parts[$i].Trim()
        }

        if ($parts.Length -ge 2) {
            $newParts = New-Object 'object[]' $parts.Length
            $newParts[0] = $parts[1]
            $newParts[1] = $parts[0]
            for ($i = 2; $i -lt $parts.Length; $i++) {
                $newParts[$i] = $parts[$i]
            }
            $cell.Value2 = [string]::Join(", ", $newParts)
        }
    }
}
